$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G10").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G11").Value = "Shivamogga (Shimoga)"
$ws.Range("G17").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G21").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G23").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G28").Value = "Shivamogga (Shimoga)"
$ws.Range("G34").Value = "Vijayapura (Bijapur)"
$ws.Range("G36").Value = "Shivamogga (Shimoga)"
$ws.Range("G45").Value = "Shivamogga (Shimoga)"
$ws.Range("G46").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G48").Value = "Shivamogga (Shimoga)"
